$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force values to be written as literal text (matches
# the source file, which stores every B:E cell as an inline string) instead
# of letting Excel auto-convert numeric-looking text such as "1.00" or
# "0.0000218" into real numbers when assigned through .Value.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"

function Set-TextValue($cellRef, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# Row 2
Set-TextValue "D2" '61.158.55'
Set-TextValue "E2" '  -4.23%  '

# Row 3
Set-TextValue "D3" '2.973.79'
Set-TextValue "E3" '  -3.78%  '

# Row 4
Set-TextValue "E4" '  -0.01%  '

# Row 5
Set-TextValue "D5" '545.01'
Set-TextValue "E5" '  +0.28%  '

# Row 6
Set-TextValue "D6" '131.91'
Set-TextValue "E6" '  -5.30%  '

# Row 7
Set-TextValue "E7" '  -0.04%  '

# Row 8
Set-TextValue "D8" '2.969.52'
Set-TextValue "E8" '  -3.75%  '

# Row 9
Set-TextValue "E9" '  -1.23%  '

# Row 10
Set-TextValue "E10" '  -6.59%  '

# Row 11
Set-TextValue "D11" '5.85'
Set-TextValue "E11" '  -10.56%  '

# Row 12
Set-TextValue "D12" '0.443'
Set-TextValue "E12" '  -2.84%  '

# Row 13
Set-TextValue "D13" '0.0000218'
Set-TextValue "E13" '  -3.20%  '

# Row 14
Set-TextValue "D14" '33.83'
Set-TextValue "E14" '  -2.75%  '

# Row 15
Set-TextValue "D15" '3.457.83'
Set-TextValue "E15" '  -3.65%  '

# Row 16
Set-TextValue "D16" '0.109'
Set-TextValue "E16" '  -2.96%  '

# Row 17
Set-TextValue "D17" '61.125.00'
Set-TextValue "E17" '  -4.40%  '

# Row 18
Set-TextValue "D18" '2.978.73'
Set-TextValue "E18" '  -3.57%  '

# Row 19
Set-TextValue "D19" '6.58'
Set-TextValue "E19" '  -1.18%  '

# Row 20
Set-TextValue "D20" '466.40'
Set-TextValue "E20" '  -2.88%  '

# Row 21
Set-TextValue "D21" '13.08'
Set-TextValue "E21" '  -2.65%  '

# Row 22
Set-TextValue "D22" '0.665'
Set-TextValue "E22" '  -4.83%  '

# Row 23
Set-TextValue "D23" '6.95'
Set-TextValue "E23" '  -1.94%  '

# Row 24
Set-TextValue "D24" '79.70'
Set-TextValue "E24" '  +1.02%  '

# Row 25
Set-TextValue "D25" '11.91'
Set-TextValue "E25" '  -3.58%  '

# Row 26
Set-TextValue "D26" '0.999'
Set-TextValue "E26" '  -0.20%  '

# Row 27
Set-TextValue "D27" '2.69'
Set-TextValue "E27" '  -1.26%  '

# Row 28
Set-TextValue "D28" '7.63'
Set-TextValue "E28" '  -5.41%  '

# Row 29
Set-TextValue "D29" '1.00'
Set-TextValue "E29" '  +0.14%  '

# Row 30
Set-TextValue "D30" '1.88'
Set-TextValue "E30" '  -1.14%  '

# Row 31
Set-TextValue "D31" '25.31'
Set-TextValue "E31" '  -3.86%  '

# Row 32
Set-TextValue "E32" '  -3.08%  '

# Row 33
Set-TextValue "D33" '2.28'
Set-TextValue "E33" '  -2.88%  '

# Row 34
Set-TextValue "B34" 'OKB'
Set-TextValue "C34" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D34" '54.98'
Set-TextValue "E34" '  -3.75%  '

# Row 35
Set-TextValue "B35" 'NEARProtocol'
Set-TextValue "C35" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D35" '5.42'
Set-TextValue "E35" '  +1.10%  '

# Row 36
Set-TextValue "D36" '5.82'
Set-TextValue "E36" '  -3.27%  '

# Row 37
Set-TextValue "D37" '442.55'
Set-TextValue "E37" '  -10.19%  '

# Row 38
Set-TextValue "D38" '3.155.15'
Set-TextValue "E38" '  -3.08%  '

# Row 39
Set-TextValue "D39" '0.0786'
Set-TextValue "E39" '  -1.61%  '

# Row 40
Set-TextValue "D40" '0.0376'
Set-TextValue "E40" '  -7.16%  '

# Row 41
Set-TextValue "E41" '  -2.77%  '

# Row 42
Set-TextValue "D42" '8.05'
Set-TextValue "E42" '  -0.77%  '

# Row 43
Set-TextValue "D43" '2.37'
Set-TextValue "E43" '  -11.28%  '

# Row 45
Set-TextValue "D45" '25.63'
Set-TextValue "E45" '  +1.89%  '

# Row 46
Set-TextValue "D46" '0.240'
Set-TextValue "E46" '  -5.29%  '

# Row 47
Set-TextValue "E47" '  -1.75%  '

# Row 48
Set-TextValue "D48" '116.75'
Set-TextValue "E48" '  -5.95%  '

# Row 49
Set-TextValue "D49" '1.94'
Set-TextValue "E49" '  -4.82%  '

# Row 50
Set-TextValue "E50" '  +7.95%  '

# Row 51
Set-TextValue "E51" '  -9.05%  '

$ws.Application.CutCopyMode = $false
$scratch.Clear()